# Update latest output (run 119)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Schedule" sheet: the optimiser now returns 2 pump-on windows instead of 3.
# Window 1 (04:00-16:00, 12h) absorbs the old 04:00-08:00 "off" gap, and the
# trailing 20:00-00:00 window's numbers move up into row 3. Row 4 is removed.
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Cells.Item(2, 1).Value = 46042.16666666666
$schedule.Cells.Item(2, 2).Value = 46042.66666666666
$schedule.Cells.Item(2, 3).Value = 12
$schedule.Cells.Item(2, 4).Value = 45.36
$schedule.Cells.Item(2, 5).Value = 336.289161
$schedule.Cells.Item(2, 6).Value = 7.413782208994709

$schedule.Cells.Item(3, 1).Value = 46042.83333333334
$schedule.Cells.Item(3, 2).Value = 46043
$schedule.Cells.Item(3, 3).Value = 4
$schedule.Cells.Item(3, 4).Value = 15.12
$schedule.Cells.Item(3, 5).Value = 439.2953565
$schedule.Cells.Item(3, 6).Value = 29.05392569444444

# Old row 4 no longer exists - drop it (shrinks dimension to A1:F3).
$schedule.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# "Detailed" sheet: refreshed price curve (col B) for this run, the
# historical/forecast cutover (col C) shifted from row 13 to row 15, and the
# pump status (col E) "OFF" window moved from rows 2-9 to rows 10-17.
# ---------------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

# Pump_Status: rows 2-9 switch ON -> OFF
for ($r = 2; $r -le 9; $r++) {
    $detailed.Cells.Item($r, 5).Value = "OFF"
}

# Pump_Status: rows 10-17 switch OFF -> ON
for ($r = 10; $r -le 17; $r++) {
    $detailed.Cells.Item($r, 5).Value = "ON"
}

# Type: rows 13-14 switch forecast -> historical
$detailed.Cells.Item(13, 3).Value = "historical"
$detailed.Cells.Item(14, 3).Value = "historical"

# Price: refreshed values for rows 10-49 (row 18 unchanged)
$priceUpdates = @{
    10 = 57.06003
    11 = 61.7706
    12 = 65
    13 = 77.52726
    14 = 77.94
    15 = 65
    16 = 36.05924
    17 = 0.7
    19 = 0.00002
    20 = -4.22824
    21 = -6
    22 = -6.62627
    23 = -6
    24 = -7.19965
    25 = -6.31573
    26 = -6.77908
    27 = -8.05312
    28 = -8.5175
    29 = -7.74764
    30 = -7.92377
    31 = -7.93011
    32 = -7.18828
    33 = -6.3358
    34 = -3.05349
    35 = -5.47652
    36 = -4.57235
    37 = 9.65405
    38 = 9.85617
    39 = 33.01118
    40 = 55.50964
    41 = 57.3
    42 = 59.09576
    43 = 57.3
    44 = 56.09936
    45 = 49.42419
    46 = 56.98
    47 = 57.3
    48 = 57.3
    49 = 57.06003
}

foreach ($r in $priceUpdates.Keys) {
    $detailed.Cells.Item($r, 2).Value = $priceUpdates[$r]
}
